# Add a new custom paragraph style "byline", based on the existing
# "Author" style, matching the style block added to styles.xml
# (and the parallel stylesWithEffects.xml) in the target revision:
#
#   <w:style w:type="paragraph" w:customStyle="1" w:styleId="byline">
#     <w:name w:val="byline"/>
#     <w:basedOn w:val="Author"/>
#     <w:qFormat/>
#     <w:rsid w:val="00076B7A"/>
#   </w:style>

$d = $word.ActiveDocument

# wdStyleTypeParagraph = 1
$byline = $d.Styles.Add("byline", 1)
$byline.BaseStyle = "Author"
$byline.QuickStyle = $true
